# Update the Feefo / Review report rows (data/recent_reviews_web_ready.xlsx sheet)
# with the refreshed 2025-04-28 product/review snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('1019839', 'Easter Chocolate Mosaic Gift With Personalised Message', 'cocoapod', 2, 'https://www.notonthehighstreet.com/cocoapod/product/easter-chocolate-mosaic-gift-with-personalised-message', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1019839&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('1067222', 'Personalised Special Year Ginger Whisky With Music', 'mixpixie', 2, 'https://www.notonthehighstreet.com/mixpixie/product/personalised-special-year-ginger-whisky', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1067222&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('1087692', '''''Happy Birthday!'' Rainbow Candles Card', 'gabriellesollyillustration', 2, 'https://www.notonthehighstreet.com/gabriellesollyillustration/product/it-s-your-birthday-rainbow-candles-card', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1087692&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('1162095', 'Personalised 1955 70th Birthday Threepence Keyring', 'charlieboots', 2, 'https://www.notonthehighstreet.com/charlieboots/product/personalised-1952-70th-birthday-threepence-keyring', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1162095&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('1263737', 'Personalised Mug ''Used To Work With Absolute Legends''', 'arrowgiftco', 2, 'https://www.notonthehighstreet.com/arrowgiftco/product/personalised-mug-used-to-work-with-absolute-legends', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1263737&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('1278588', 'Bunny In A Box', 'thegourmetchocolatepizzaco', 3, 'https://www.notonthehighstreet.com/thegourmetchocolatepizzaco/product/bunny-in-a-box', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1278588&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('1285611', 'Personalised My First Easter Bunny Toy With Name', 'looppa', 3, 'https://www.notonthehighstreet.com/looppa/product/personalised-my-first-easter-bunny-toy-with-name', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1285611&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('1353406', 'F1 Formula One 2025 Calendar Track T Shirt Gift For Him', 'nappyhead', 2, 'https://www.notonthehighstreet.com/nappyhead/product/formula-1-track-t-shirt-gift-for-him', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1353406&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('1383730', 'Personalised Relax And Unwind Spa Pamper Hamper', 'mixpixie', 2, 'https://www.notonthehighstreet.com/mixpixie/product/personalised-relax-unwind-spa-pamper-hamper', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1383730&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('1435766', 'Sterling Silver Enamel Bluebell Flower Drop Hook Earrings', 'silverrainsilver', 2, 'https://www.notonthehighstreet.com/silverrainsilver/product/sterling-silver-enamel-bluebell-flower-drop-hook-earrings', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1435766&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('1470941', 'Porcelain Bird Earrings French Vintage Charm', 'uniques', 2, 'https://www.notonthehighstreet.com/uniques/product/porcelain-bird-earrings-french-vintage-charm', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1470941&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('1488236', 'Personalized Jewelry Box', 'lineamasks', 2, 'https://www.notonthehighstreet.com/lineamasks/product/personalized-jewelry-box', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1488236&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('1489841', 'Simnel Cake', 'theoriginalcakecompany', 2, 'https://www.notonthehighstreet.com/theoriginalcakecompany/product/simnel-cake', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1489841&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('1491712', 'Statement Maker Socks Box - Letterbox Gift, Three Pairs', 'standoutsocks', 2, 'https://www.notonthehighstreet.com/standoutsocks/product/statement-maker-socks-box-letterbox-gift-three-pairs', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1491712&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('159976', 'Personalised Pet Pebble With A Pawprint', 'letterfestengraving', 2, 'https://www.notonthehighstreet.com/letterfestengraving/product/personalised-pet-pebble-with-a-pawprint', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=159976&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('298443', 'Personalised Metal Street Sign', 'oakdenedesigns', 2, 'https://www.notonthehighstreet.com/oakdenedesigns/product/personalised-metal-street-sign', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=298443&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('469358', 'Family Birthstone Link Bracelet', 'joybycorrinesmith', 3, 'https://www.notonthehighstreet.com/joybycorrinesmith/product/family-birthstone-link-bracelet', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=469358&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('748820', 'Personalised Metallic Sound Wave Song Print', 'mixpixie', 2, 'https://www.notonthehighstreet.com/mixpixie/product/personalised-metallic-foil-sound-wave-song-print', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=748820&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('771041', 'Personalised Relaxation And Wellness Pamper Gift Set', 'bankslyonbotanical', 2, 'https://www.notonthehighstreet.com/bankslyonbotanical/product/relaxation-gift-set', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=771041&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('786481', 'Sterling Silver Birthstone Sliding Bracelet', 'gaamaa', 2, 'https://www.notonthehighstreet.com/gaamaa/product/sterling-silver-birthstone-sliding-bracelet', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=786481&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('808656', 'Happy Easter Daffodil Forever Flowers Box', 'noordinarygiftcompany', 2, 'https://www.notonthehighstreet.com/noordinarygiftcompany/product/forever-daffodils-keepsake', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=808656&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('876141', 'Personalised Birth Flower Scarf', 'theforestandco', 2, 'https://www.notonthehighstreet.com/theforestandco/product/personalised-birth-flower-scarf', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=876141&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('887531', 'Personalised Couples Name Heart Cushion Wedding Gift', 'thatsnicethat', 2, 'https://www.notonthehighstreet.com/thatsnicethat/product/personalised-couples-name-cushion', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=887531&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('915034', 'Personalised Wedding Day / First Anniversary Card', 'sascreative', 2, 'https://www.notonthehighstreet.com/sascreative/product/personalised-wedding-day-first-anniversary-card', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=915034&displayFeedbackType=PRODUCT&timeFrame=ALL'),
    @('932102', 'Personalised Meaningful Message Leather Bookmark', 'lisaangeljewellery', 2, 'https://www.notonthehighstreet.com/lisaangeljewellery/product/personalised-meaningful-message-leather-bookmark', 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=932102&displayFeedbackType=PRODUCT&timeFrame=ALL'),
)

# Product Code values are numeric-looking strings (e.g. "1019839") that must stay
# plain text, matching the workbook's original inlineStr encoding - force the
# column to a text number format before writing so COM doesn't coerce them to
# numbers.
$codeRange = $ws.Range("A2:A" + (1 + $data.Count))
$codeRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $titleCell = $ws.Cells.Item($r, 2)
    $titleCell.Value = $row[1]
    # Titles whose literal text begins with a single quote (') get the leading
    # character swallowed by the COM layer's "quote prefix" auto-formatting and
    # get flagged with a quotePrefix style. Writing the value with a doubled
    # leading quote survives that stripping and leaves exactly one quote behind,
    # then ClearFormats drops the stray quotePrefix style so the cell matches a
    # plain, unstyled text cell again.
    if ($row[1].StartsWith("'")) {
        $titleCell.ClearFormats()
    }
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
